# Updated cryptos list on Fri Jul 28 18:42:13 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row,
# and fixes the EnergySwap / BabyDogeCoin row ordering (rows 49-50).
#
# Note: several Price values look numeric (e.g. "1.0000", "18.50",
# "0.000008260") but must stay verbatim text, matching the source feed's
# formatting. A leading apostrophe forces Excel to store them as text
# instead of re-parsing/normalizing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.281.62'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.873.44'
$ws.Range('E3').Value = '  +0.63%  '
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').Value = '''0.7121'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').Value = '''241.63'
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('D7').Value = '''1.0000'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').Value = '''0.3104'
$ws.Range('E8').Value = '  +1.23%  '
$ws.Range('D9').Value = '''0.07712'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').Value = '''25.06'
$ws.Range('E10').Value = '  +0.72%  '
$ws.Range('D11').Value = '''0.08381'
$ws.Range('E11').Value = '  +1.70%  '
$ws.Range('D12').Value = '1.896.77'
$ws.Range('E12').Value = '  +1.62%  '
$ws.Range('D13').Value = '''5.216'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('E14').Value = '  -0.57%  '
$ws.Range('D15').Value = '''91.27'
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('D16').Value = '29.295.71'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = '''0.000008260'
$ws.Range('E17').Value = '  +6.40%  '
$ws.Range('D18').Value = '''5.974'
$ws.Range('E18').Value = '  +2.57%  '
$ws.Range('D19').Value = '''242.35'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').Value = '2.129.34'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').Value = '''13.20'
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').Value = '''7.833'
$ws.Range('E23').Value = '  -1.16%  '
$ws.Range('E24').Value = '  -0.38%  '
$ws.Range('D25').Value = '''0.1619'
$ws.Range('E25').Value = '  +2.31%  '
$ws.Range('D26').Value = '''163.27'
$ws.Range('E26').Value = '  +0.90%  '
$ws.Range('D27').Value = '''9.013'
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('D28').Value = '''18.50'
$ws.Range('E28').Value = '  +2.23%  '
$ws.Range('E29').Value = '  +0.86%  '
$ws.Range('E30').Value = '  +1.94%  '
$ws.Range('E31').Value = '  -3.76%  '
$ws.Range('D32').Value = '''4.320'
$ws.Range('E32').Value = '  +6.19%  '
$ws.Range('E33').Value = '  +1.42%  '
$ws.Range('D34').Value = '''1.925'
$ws.Range('E34').Value = '  +0.78%  '
$ws.Range('D35').Value = '''0.7484'
$ws.Range('E35').Value = '  +2.96%  '
$ws.Range('E36').Value = '  -0.17%  '
$ws.Range('D37').Value = '''2.681'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').Value = '''0.01857'
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('D39').Value = '''2.715'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('D40').Value = '1.156.47'
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('D41').Value = '''6.370'
$ws.Range('E41').Value = '  +4.69%  '
$ws.Range('D42').Value = '''73.11'
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('D43').Value = '''0.8845'
$ws.Range('E43').Value = '  -1.64%  '
$ws.Range('D44').Value = '''105.29'
$ws.Range('E44').Value = '  +3.57%  '
$ws.Range('D45').Value = '''0.9992'
$ws.Range('E45').Value = '  -0.30%  '
$ws.Range('D46').Value = '2.027.18'
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('D47').Value = '''1.803'
$ws.Range('E47').Value = '  +2.58%  '
$ws.Range('D48').Value = '''0.5191'
$ws.Range('E48').Value = '  -1.90%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '''0.00000000120'
$ws.Range('E49').Value = '  +2.85%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''9.378'
$ws.Range('E50').Value = '  +1.75%  '
$ws.Range('E51').Value = '  +1.85%  '
